# "adding low level flowcharts"
# Mark four previously-open audit comments (LLD/flowchart related rows)
# as closed, and record who they were assigned to / closed by.
#
#   F4  / G4  : comment_03 (CM)     open -> closed, Assigned to -> Osama
#   F8  / G8  : comment_07 (reviews) open -> closed, Assigned to -> Medhat
#   F13 / G13 : comment_12 (Design) open -> closed, Assigned to -> Nada
#   F14 / G14 : comment_13 (Design) open -> closed, Assigned to -> Nada

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("state" column F2) already carries the exact "closed" formatting
# (font/style) that these rows need to adopt, so copy it onto each target
# cell before giving it the "closed" text.
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("G4").Value = "Osama"

$ws.Range("F2").Copy($ws.Range("F8"))
$ws.Range("G8").Value = "Medhat"

$ws.Range("F2").Copy($ws.Range("F13"))
$ws.Range("G13").Value = "Nada"

$ws.Range("F2").Copy($ws.Range("F14"))
$ws.Range("G14").Value = "Nada"

# Move the active selection to F8, matching the saved workbook state.
$ws.Range("F8").Select()
